$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (ALC) - hunk 0
$ws.Range("H6").Value = 1294.9
$ws.Range("I6").Value = 1405.2222
$ws.Range("J6").Value = 302
$ws.Range("K6").Value = 4215.6666
$ws.Range("L6").Value = 906
$ws.Range("M6").Value = -4103.6666
$ws.Range("N6").Value = -1130

# Row 33 (ALC) - hunk 1
$ws.Range("H33").Value = 388.92307
$ws.Range("I33").Value = 147.11765
$ws.Range("K33").Value = 147.11765
$ws.Range("M33").Value = 81.88235

# Row 70 (ALC) - hunk 2
$ws.Range("H70").Value = 50001228
$ws.Range("I70").Value = 1495
$ws.Range("J70").Value = 62501164
$ws.Range("K70").Value = 4485
$ws.Range("L70").Value = 187503492
$ws.Range("M70").Value = -4215
$ws.Range("N70").Value = -187504032

# Row 73 (ALC) - hunk 3
$ws.Range("H73").Value = 50001228
$ws.Range("I73").Value = 1495
$ws.Range("J73").Value = 62501164
$ws.Range("K73").Value = 4485
$ws.Range("L73").Value = 187503492
$ws.Range("M73").Value = -3549
$ws.Range("N73").Value = -187505364

# Row 113 (ALC) - hunk 4
$ws.Range("H113").Value = 4072
$ws.Range("I113").Value = 3149
$ws.Range("J113").Value = 4892.4443
$ws.Range("K113").Value = 3149
$ws.Range("L113").Value = 4892.4443
$ws.Range("M113").Value = 105
$ws.Range("N113").Value = -11400.4443

# Row 118 (ALC) - hunk 5
$ws.Range("H118").Value = 1003.5714
$ws.Range("I118").Value = 873
$ws.Range("J118").Value = 1101.5
$ws.Range("K118").Value = 2619
$ws.Range("L118").Value = 3304.5
$ws.Range("M118").Value = -962
$ws.Range("N118").Value = -6618.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM) - hunk 6
$ws.Range("H61").Value = 3507.7273
$ws.Range("I61").Value = 2933.7222
$ws.Range("K61").Value = 2933.7222
$ws.Range("M61").Value = -2721.7222

# Row 122 (ARM) - hunk 7
$ws.Range("H122").Value = 4153.775
$ws.Range("I122").Value = 3564.75
$ws.Range("J122").Value = 5037.3125
$ws.Range("K122").Value = 10694.25
$ws.Range("L122").Value = 15111.9375
$ws.Range("M122").Value = -8244.25
$ws.Range("N122").Value = -20011.9375

# Row 132 (ARM) - hunk 8
$ws.Range("H132").Value = 2448.5833
$ws.Range("I132").Value = 2238.85
$ws.Range("K132").Value = 6716.549999999999
$ws.Range("M132").Value = -4186.549999999999

# Row 136 (ARM) - hunk 9
$ws.Range("H136").Value = 3507.7273
$ws.Range("I136").Value = 2933.7222
$ws.Range("K136").Value = 8801.1666
$ws.Range("M136").Value = -6251.1666

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM) - hunk 10
$ws.Range("H20").Value = 22529.666
$ws.Range("I20").Value = 37178
$ws.Range("J20").Value = 2022
$ws.Range("K20").Value = 37178
$ws.Range("L20").Value = 2022
$ws.Range("M20").Value = -36931
$ws.Range("N20").Value = -2516

# Row 54 (BSM) - hunk 11
$ws.Range("H54").Value = 2752.4443
$ws.Range("I54").Value = 2397.8333
$ws.Range("K54").Value = 2397.8333
$ws.Range("M54").Value = -1913.8333

# Row 86 (BSM) - hunk 12
$ws.Range("H86").Value = 5399.222
$ws.Range("J86").Value = 5571.2856
$ws.Range("L86").Value = 5571.2856
$ws.Range("N86").Value = -7817.2856

# Row 89 (BSM) - hunk 13
$ws.Range("H89").Value = 5399.222
$ws.Range("J89").Value = 5571.2856
$ws.Range("L89").Value = 27856.428
$ws.Range("N89").Value = -39088.428

# Row 94 (BSM) - hunk 14
$ws.Range("H94").Value = 627.04346
$ws.Range("I94").Value = 675.381
$ws.Range("J94").Value = 119.5
$ws.Range("K94").Value = 675.381
$ws.Range("L94").Value = 119.5
$ws.Range("M94").Value = -224.381
$ws.Range("N94").Value = -1021.5

# Row 134 (BSM) - hunk 15
$ws.Range("H134").Value = 2472509
$ws.Range("I134").Value = 2567524.8
$ws.Range("K134").Value = 7702574.399999999
$ws.Range("M134").Value = -7700039.399999999

$ws = $wb.Worksheets.Item("CRP")
# Row 15 (CRP) - hunk 16
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

# Row 22 (CRP) - hunk 17
$ws.Range("H22").Value = 2224.6924
$ws.Range("I22").Value = 2422.1
$ws.Range("J22").Value = 1566.6666
$ws.Range("K22").Value = 2422.1
$ws.Range("L22").Value = 1566.6666
$ws.Range("M22").Value = -2072.1
$ws.Range("N22").Value = -2266.6666

# Row 31 (CRP) - hunk 18
$ws.Range("H31").Value = 4065.8125
$ws.Range("I31").Value = 1239.8334
$ws.Range("K31").Value = 1239.8334
$ws.Range("M31").Value = -944.8334

# Row 34 (CRP) - hunk 19
$ws.Range("H34").Value = 4065.8125
$ws.Range("I34").Value = 1239.8334
$ws.Range("K34").Value = 1239.8334
$ws.Range("M34").Value = -1037.8334

# Row 134 (CRP) - hunk 20
$ws.Range("H134").Value = 1424.3334
$ws.Range("I134").Value = 1462.4546
$ws.Range("J134").Value = 1005
$ws.Range("K134").Value = 4387.3638
$ws.Range("L134").Value = 3015
$ws.Range("M134").Value = -1852.3638
$ws.Range("N134").Value = -8085

$ws = $wb.Worksheets.Item("CUL")
# Row 6 (CUL) - hunk 21
$ws.Range("H6").Value = 22.571428
$ws.Range("I6").Value = 26
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 78
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 35
$ws.Range("N6").Value = -232

$ws = $wb.Worksheets.Item("GSM")
# Row 103 (GSM) - hunk 22
$ws.Range("H103").Value = 91530.8
$ws.Range("J103").Value = 91530.8
$ws.Range("L103").Value = 91530.8
$ws.Range("N103").Value = -93874.8

# Row 113 (GSM) - hunk 23
$ws.Range("H113").Value = 9227.030000000001
$ws.Range("I113").Value = 6886.8965
$ws.Range("J113").Value = 26193
$ws.Range("K113").Value = 6886.8965
$ws.Range("L113").Value = 26193
$ws.Range("M113").Value = -4716.8965
$ws.Range("N113").Value = -30533

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW) - hunk 24
$ws.Range("H61").Value = 3767.1667
$ws.Range("I61").Value = 1602.5834
$ws.Range("J61").Value = 5931.75
$ws.Range("K61").Value = 1602.5834
$ws.Range("L61").Value = 5931.75
$ws.Range("M61").Value = -1400.5834
$ws.Range("N61").Value = -6335.75

# Row 68 (LTW) - hunk 25
$ws.Range("H68").Value = 2140.8333
$ws.Range("I68").Value = 2062.4546
$ws.Range("J68").Value = 3003
$ws.Range("K68").Value = 2062.4546
$ws.Range("L68").Value = 3003
$ws.Range("M68").Value = -1313.4546
$ws.Range("N68").Value = -4501

# Row 71 (LTW) - hunk 26
$ws.Range("H71").Value = 2140.8333
$ws.Range("I71").Value = 2062.4546
$ws.Range("J71").Value = 3003
$ws.Range("K71").Value = 10312.273
$ws.Range("L71").Value = 15015
$ws.Range("M71").Value = -6568.273000000001
$ws.Range("N71").Value = -22503

# Row 93 (LTW) - hunk 27
$ws.Range("H93").Value = 24391740
$ws.Range("I93").Value = 40001260
$ws.Range("J93").Value = 1860.75
$ws.Range("K93").Value = 40001260
$ws.Range("L93").Value = 1860.75
$ws.Range("M93").Value = -40000012
$ws.Range("N93").Value = -4356.75

# Row 102 (LTW) - hunk 28
$ws.Range("H102").Value = 121998.5
$ws.Range("J102").Value = 121998.5
$ws.Range("L102").Value = 121998.5
$ws.Range("N102").Value = -128488.5

# Row 113 (LTW) - hunk 29
$ws.Range("H113").Value = 3767.1667
$ws.Range("I113").Value = 1602.5834
$ws.Range("J113").Value = 5931.75
$ws.Range("K113").Value = 1602.5834
$ws.Range("L113").Value = 5931.75
$ws.Range("M113").Value = 567.4166
$ws.Range("N113").Value = -10271.75

# Row 119 (LTW) - hunk 30
$ws.Range("H119").Value = 99913
$ws.Range("J119").Value = 99913
$ws.Range("L119").Value = 99913
$ws.Range("N119").Value = -109589

$ws = $wb.Worksheets.Item("WVR")
# Row 102 (WVR) - hunk 31
$ws.Range("H102").Value = 74883.5
$ws.Range("J102").Value = 74883.5
$ws.Range("L102").Value = 74883.5
$ws.Range("N102").Value = -81373.5

# Row 132 (WVR) - hunk 32
$ws.Range("H132").Value = 2205.4138
$ws.Range("I132").Value = 1904.1904
$ws.Range("K132").Value = 5712.5712
$ws.Range("M132").Value = -3182.5712
